$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column cells are treated as plain text, matching the
# source data (which stores prices like "35.208.08" or "8.00" as text,
# not numbers), so trailing zeros / multi-dot groupings are preserved.
$priceCells = "D2","D3","D5","D8","D11","D12","D14","D17","D19","D20","D21","D25","D26","D27","D31","D33","D36","D37","D40","D41","D42","D43","D46","D49"
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "35.208.08"
$ws.Range("E2").Value = "  +1.18%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.859.20"
$ws.Range("E3").Value = "  +1.66%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.41%  "

# Row 5 - BNB
$ws.Range("D5").Value = "239.43"
$ws.Range("E5").Value = "  +3.74%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.80%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.36%  "

# Row 8 - Solana
$ws.Range("D8").Value = "42.12"
$ws.Range("E8").Value = "  +7.35%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +0.85%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.44%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.0989"
$ws.Range("E11").Value = "  +0.02%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.128.46"
$ws.Range("E12").Value = "  +1.64%  "

# Row 13 - Chainlink
$ws.Range("E13").Value = "  +1.62%  "

# Row 14 - WrappedEther
$ws.Range("D14").Value = "1.862.56"
$ws.Range("E14").Value = "  +1.90%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +1.17%  "

# Row 16 - Polkadot
$ws.Range("E16").Value = "  +2.01%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "35.167.73"
$ws.Range("E17").Value = "  +1.16%  "

# Row 18 - Litecoin
$ws.Range("E18").Value = "  +0.63%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "0.0₃0797"
$ws.Range("E19").Value = "  +1.54%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "240.56"
$ws.Range("E20").Value = "  +0.24%  "

# Row 21 - Avalanche
$ws.Range("D21").Value = "12.22"
$ws.Range("E21").Value = "  +0.59%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  +1.97%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.50%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +0.83%  "

# Row 25 - Monero
$ws.Range("D25").Value = "168.70"
$ws.Range("E25").Value = "  -1.84%  "

# Row 26 - PancakeSwap
$ws.Range("D26").Value = "1.91"
$ws.Range("E26").Value = "  +27.48%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "8.00"
$ws.Range("E27").Value = "  +3.48%  "

# Row 28 - EthereumClassic
$ws.Range("E28").Value = "  +1.85%  "

# Row 29 - Stellar
$ws.Range("E29").Value = "  +0.30%  "

# Row 30 - BinanceUSD
$ws.Range("E30").Value = "  +0.41%  "

# Row 31 - Hedera
$ws.Range("D31").Value = "0.0559"
$ws.Range("E31").Value = "  +1.62%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +2.39%  "

# Row 33 - WEMIXToken
$ws.Range("D33").Value = "1.82"
$ws.Range("E33").Value = "  +27.64%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  +2.26%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  +10.49%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "0.814"
$ws.Range("E36").Value = "  +16.66%  "

# Row 37 - TrustWalletToken
$ws.Range("D37").Value = "1.32"
$ws.Range("E37").Value = "  +7.72%  "

# Row 38 - ARBITRUM
$ws.Range("E38").Value = "  +4.78%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  +4.06%  "

# Row 40 - Aave
$ws.Range("D40").Value = "89.93"
$ws.Range("E40").Value = "  -1.42%  "

# Row 41 - Maker
$ws.Range("D41").Value = "1.346.50"
$ws.Range("E41").Value = "  +0.37%  "

# Row 42 - now InjectiveProtocol (was Kaspa)
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "14.96"
$ws.Range("E42").Value = "  +3.30%  "

# Row 43 - now Kaspa (was InjectiveProtocol)
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "0.0588"
$ws.Range("E43").Value = "  +13.16%  "

# Row 44 - RenderToken
$ws.Range("E44").Value = "  +3.51%  "

# Row 45 - HuobiToken
$ws.Range("E45").Value = "  +0.19%  "

# Row 46 - Gas
$ws.Range("D46").Value = "12.30"
$ws.Range("E46").Value = "  +41.74%  "

# Row 48 - FraxShare
$ws.Range("E48").Value = "  +5.41%  "

# Row 49 - RocketPoolETH
$ws.Range("D49").Value = "2.047.00"
$ws.Range("E49").Value = "  +1.83%  "

# Row 50 - Cronos
$ws.Range("E50").Value = "  +1.11%  "

# Row 51 - PaxDollar
$ws.Range("E51").Value = "  +0.39%  "
